# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Rule "R40" (row 11) is renamed to "1"
$ws.Range("B11").Value = "1"
